$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as literal text even though it
# looks like a number/boolean (e.g. "8", "1923", "true", "47200"), matching
# the source workbook where these are plain shared-string entries (t="s")
# rather than numeric/boolean cells. A leading apostrophe forces text entry;
# ClearFormats() afterwards drops the transient "quote prefix" cell format
# that the apostrophe trick applies, so the cell keeps the plain default
# style instead of picking up a new one.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).ClearFormats()
}

# Fill in the remaining columns (F:S) for the existing last data row (row 8)
Set-TextValue "F8" "8"
Set-TextValue "G8" "1923"
Set-TextValue "H8" "true"
Set-TextValue "I8" "true"
$ws.Range("J8").Value = "Luke"
$ws.Range("K8").Value = "Botsford"
$ws.Range("L8").Value = "Carter-Sporer"
$ws.Range("M8").Value = "4820 Foster Way"
$ws.Range("N8").Value = "Suite 808"
$ws.Range("O8").Value = "Israel"
$ws.Range("P8").Value = "Nebraska"
$ws.Range("Q8").Value = "Port Quinn"
Set-TextValue "R8" "47200"
$ws.Range("S8").Value = "807-647-6977"

# Add a new data row (row 9) with a new signup record
$ws.Range("A9").Value = "Mrs"
$ws.Range("B9").Value = "Luke"
$ws.Range("C9").Value = "jaclyn.brakus@yahoo.com"
$ws.Range("D9").Value = "jvwawaqlmm75td"
Set-TextValue "E9" "19"
